$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Planning poker point values added to the "Planning Poker Ratings" columns
# (C = Sam Pugh, D = Ben Kownacki, per the B4/B5 author headers) for every
# user-story row.
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 21

$ws.Range("C5").Value = 34
$ws.Range("D5").Value = 34

$ws.Range("C6").Value = 34
$ws.Range("D6").Value = 34

$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 34

$ws.Range("C8").Value = 21
$ws.Range("D8").Value = 21

$ws.Range("C9").Value = 34
$ws.Range("D9").Value = 21

$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 89

$ws.Range("C11").Value = 34
$ws.Range("D11").Value = 55

$ws.Range("C12").Value = 55
$ws.Range("D12").Value = 21

$ws.Range("C13").Value = 55
$ws.Range("D13").Value = 55

# Leave the cursor where the author's last edit landed.
$ws.Range("C13").Select()
